$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: take over former row 5 values (Especial, higher prices, later date)
$ws.Range("D2").Value = 44460
$ws.Range("L2").Value = "Especial"
$ws.Range("N2").Value = 31000
$ws.Range("O2").Value = 32000
$ws.Range("P2").Value = 31500
$ws.Range("S2").Value = 3150

# Row 3: take over former row 6 values (Primera, volumen 30, prices 30000)
$ws.Range("D3").Value = 44460
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 30000
$ws.Range("O3").Value = 30000
$ws.Range("P3").Value = 30000
$ws.Range("S3").Value = 3000

# Row 4: take over former row 2 values (date only changes)
$ws.Range("D4").Value = 44446

# Row 5: take over former row 7 values (Primera, base prices, earlier date)
$ws.Range("D5").Value = 44447
$ws.Range("L5").Value = "Primera"
$ws.Range("N5").Value = 21000
$ws.Range("O5").Value = 22000
$ws.Range("P5").Value = 21500
$ws.Range("S5").Value = 2150

# Row 6: take over former row 4 values
$ws.Range("D6").Value = 44448
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 21000
$ws.Range("O6").Value = 22000
$ws.Range("P6").Value = 21500
$ws.Range("S6").Value = 2150

# Row 7: take over former row 3 values
$ws.Range("D7").Value = 44452
